# Update the "Follows" links for a few titles, and fill in the two rows
# that were still missing their "Follows" value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Betrayer (row 26) now follows "Know No Fear" instead of "The First Heretic"
$ws.Range("D26").Value = "Know No Fear"

# Garro : Vow of Faith (row 45) previously had no "Follows" value - give it
# one, matching the left-aligned style already used by the other filled-in
# cells in this column (copy formatting from a neighboring filled cell first).
$ws.Range("D42").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("D45").Value = "The Flight of the Eisenstein,Nemesis"

# The Unremembered Empire (row 31) now also follows "Betrayer"
$ws.Range("D31").Value = "Betrayer,Battle for the Abyss,Vulkan Lives,Fallen Angels"

# Path of Heaven (row 49) previously had no "Follows" value either
$ws.Range("D42").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D49").Value = "Scars"

$excel.CutCopyMode = $false

# Leave the selection where the author ended up
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 36
